# Add "neologism" results row (row 14: columns A, B, C) and update the
# sheet view's selection/scroll position, per the commit
# "Add neologism results and graphs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: label, pass flag (styled like the other "pass" cells),
# and the p-value.
$ws.Range("A14").Value = "neologism"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.0044

# Apply the same "Good" cell style used by other pass-flag cells (e.g. B13)
# to the new pass-flag cell B14.
$ws.Range("B14").Style = $ws.Range("B13").Style

# Update the active selection / view to match the new focus on row 14.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F14").Select()
